$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.375614404678345
$ws.Range("B1").Value = 2.189620494842529
$ws.Range("C1").Value = 2.084326267242432
$ws.Range("D1").Value = 2.715598106384277
$ws.Range("E1").Value = 3.98919939994812
